# Nuevas pruebas de rendimiento
# - Corrige el texto de cabecera "Mínimo Ruputra" -> "Mínimo Ruptura"
# - Añade/actualiza los valores de rendimiento (columna B) para cada historia de usuario
# - Resalta la columna A (filas de datos) con un nuevo color de relleno
# - Actualiza la celda seleccionada a B5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix header typo in C1 ---
$ws.Range("C1").Value = "Mínimo Ruptura"

# --- Row data: user story name (unchanged) + new/updated measured value ---
$ws.Range("B2").Value = 1500
$ws.Range("B3").Value = 1500
$ws.Range("B4").Value = 1500
$ws.Range("B5").Value = 45
$ws.Range("B6").Value = 45
$ws.Range("B7").Value = 2800
$ws.Range("B8").Value = 3000
$ws.Range("B9").Value = 4800
$ws.Range("B10").Value = 3500
$ws.Range("B11").Value = 2600

# B7:B11 are brand new cells (previously empty) - copy the existing "value" cell
# style (used by B2:B6) onto them so they share the same look (fillId 2 style).
$ws.Range("B2").Copy()
$ws.Range("B7:B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Highlight column A of the data rows with a new fill color ---
# (Gold, Accent 4, Lighter 80% - theme color 7 / tint 0.79998168889431442 => RGB FFF2CC)
$dataRows = $ws.Range("A2:A11")
$dataRows.Interior.Color = 13431551

# --- Update the active selection to B5 ---
[void]$ws.Range("B5").Select()

Write-Host "edit complete"
